# Scheduled-runner market data refresh.
# Updates currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns (H:N) for a
# set of Leve rows across several job sheets with freshly pulled prices.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 462.7143
$ws.Cells.Item(17, 10).Value = 462.7143
$ws.Cells.Item(17, 12).Value = 1388.1429
$ws.Cells.Item(17, 14).Value = -1724.1429

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 3550.4
$ws.Cells.Item(19, 9).Value = 4000
$ws.Cells.Item(19, 10).Value = 3438
$ws.Cells.Item(19, 11).Value = 4000
$ws.Cells.Item(19, 12).Value = 3438
$ws.Cells.Item(19, 13).Value = -3825
$ws.Cells.Item(19, 14).Value = -3788

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 3277.575
$ws.Cells.Item(40, 9).Value = 2619.238
$ws.Cells.Item(40, 10).Value = 4005.2104
$ws.Cells.Item(40, 11).Value = 2619.238
$ws.Cells.Item(40, 12).Value = 4005.2104
$ws.Cells.Item(40, 13).Value = -2444.238
$ws.Cells.Item(40, 14).Value = -4355.2104

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 1922.3334
$ws.Cells.Item(70, 9).Value = 1139
$ws.Cells.Item(70, 10).Value = 2481.8572
$ws.Cells.Item(70, 11).Value = 3417
$ws.Cells.Item(70, 12).Value = 7445.571599999999
$ws.Cells.Item(70, 13).Value = -3147
$ws.Cells.Item(70, 14).Value = -7985.571599999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(73, 8).Value = 1922.3334
$ws.Cells.Item(73, 9).Value = 1139
$ws.Cells.Item(73, 10).Value = 2481.8572
$ws.Cells.Item(73, 11).Value = 3417
$ws.Cells.Item(73, 12).Value = 7445.571599999999
$ws.Cells.Item(73, 13).Value = -2481
$ws.Cells.Item(73, 14).Value = -9317.571599999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(111, 8).Value = 2974.75
$ws.Cells.Item(111, 10).Value = 5880
$ws.Cells.Item(111, 12).Value = 17640
$ws.Cells.Item(111, 14).Value = -23774

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(125, 8).Value = 2873.3333
$ws.Cells.Item(125, 10).Value = 2932.5
$ws.Cells.Item(125, 12).Value = 26392.5
$ws.Cells.Item(125, 14).Value = -31312.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 10448.2
$ws.Cells.Item(137, 9).Value = 2401
$ws.Cells.Item(137, 10).Value = 15813
$ws.Cells.Item(137, 11).Value = 7203
$ws.Cells.Item(137, 12).Value = 47439
$ws.Cells.Item(137, 13).Value = -4653
$ws.Cells.Item(137, 14).Value = -52539

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 3233.25
$ws.Cells.Item(138, 10).Value = 3452.5
$ws.Cells.Item(138, 12).Value = 10357.5
$ws.Cells.Item(138, 14).Value = -20637.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 11632915
$ws.Cells.Item(32, 9).Value = 12200369
$ws.Cells.Item(32, 11).Value = 12200369
$ws.Cells.Item(32, 13).Value = -12200082

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 2352.9333
$ws.Cells.Item(45, 9).Value = 1961
$ws.Cells.Item(45, 11).Value = 1961
$ws.Cells.Item(45, 13).Value = -1584

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 31317668
$ws.Cells.Item(61, 9).Value = 62501690
$ws.Cells.Item(61, 10).Value = 133650
$ws.Cells.Item(61, 11).Value = 62501690
$ws.Cells.Item(61, 12).Value = 133650
$ws.Cells.Item(61, 13).Value = -62501478
$ws.Cells.Item(61, 14).Value = -134074

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 31317668
$ws.Cells.Item(136, 9).Value = 62501690
$ws.Cells.Item(136, 10).Value = 133650
$ws.Cells.Item(136, 11).Value = 187505070
$ws.Cells.Item(136, 12).Value = 400950
$ws.Cells.Item(136, 13).Value = -187502520
$ws.Cells.Item(136, 14).Value = -406050

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 917
$ws.Cells.Item(22, 9).Value = 917
$ws.Cells.Item(22, 11).Value = 917
$ws.Cells.Item(22, 13).Value = -744

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(38, 8).Value = 103000
$ws.Cells.Item(38, 10).Value = 103000
$ws.Cells.Item(38, 12).Value = 103000
$ws.Cells.Item(38, 14).Value = -103832

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1500.3334
$ws.Cells.Item(107, 9).Value = 1268.8462
$ws.Cells.Item(107, 11).Value = 1268.8462
$ws.Cells.Item(107, 13).Value = 651.1538

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 52647.22
$ws.Cells.Item(134, 9).Value = 7966.778
$ws.Cells.Item(134, 10).Value = 213496.8
$ws.Cells.Item(134, 11).Value = 23900.334
$ws.Cells.Item(134, 12).Value = 640490.3999999999
$ws.Cells.Item(134, 13).Value = -21365.334
$ws.Cells.Item(134, 14).Value = -645560.3999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 871788
$ws.Cells.Item(31, 9).Value = 18293.223
$ws.Cells.Item(31, 10).Value = 1298535.4
$ws.Cells.Item(31, 11).Value = 18293.223
$ws.Cells.Item(31, 12).Value = 1298535.4
$ws.Cells.Item(31, 13).Value = -17998.223
$ws.Cells.Item(31, 14).Value = -1299125.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 871788
$ws.Cells.Item(34, 9).Value = 18293.223
$ws.Cells.Item(34, 10).Value = 1298535.4
$ws.Cells.Item(34, 11).Value = 18293.223
$ws.Cells.Item(34, 12).Value = 1298535.4
$ws.Cells.Item(34, 13).Value = -18091.223
$ws.Cells.Item(34, 14).Value = -1298939.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 5024.8
$ws.Cells.Item(3, 9).Value = 2109.8
$ws.Cells.Item(3, 10).Value = 7939.8
$ws.Cells.Item(3, 11).Value = 6329.400000000001
$ws.Cells.Item(3, 12).Value = 23819.4
$ws.Cells.Item(3, 13).Value = -6217.400000000001
$ws.Cells.Item(3, 14).Value = -24043.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(56, 8).Value = 10990.4
$ws.Cells.Item(56, 9).Value = 10990.4
$ws.Cells.Item(56, 11).Value = 10990.4
$ws.Cells.Item(56, 13).Value = -10460.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(75, 8).Value = 121354370
$ws.Cells.Item(75, 10).Value = 41861944
$ws.Cells.Item(75, 12).Value = 125585832
$ws.Cells.Item(75, 14).Value = -125587828

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(78, 8).Value = 121354370
$ws.Cells.Item(78, 10).Value = 41861944
$ws.Cells.Item(78, 12).Value = 376757496
$ws.Cells.Item(78, 14).Value = -376767480

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(82, 8).Value = 8232.571
$ws.Cells.Item(82, 10).Value = 9997.5
$ws.Cells.Item(82, 12).Value = 29992.5
$ws.Cells.Item(82, 14).Value = -30804.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(85, 8).Value = 8232.571
$ws.Cells.Item(85, 10).Value = 9997.5
$ws.Cells.Item(85, 12).Value = 29992.5
$ws.Cells.Item(85, 14).Value = -32800.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(138, 8).Value = 1962.3334
$ws.Cells.Item(138, 9).Value = 1994.8
$ws.Cells.Item(138, 10).Value = 1800
$ws.Cells.Item(138, 11).Value = 5984.4
$ws.Cells.Item(138, 12).Value = 5400
$ws.Cells.Item(138, 13).Value = -844.3999999999996
$ws.Cells.Item(138, 14).Value = -15680

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(124, 8).Value = 109978.6
$ws.Cells.Item(124, 10).Value = 109978.6
$ws.Cells.Item(124, 12).Value = 109978.6
$ws.Cells.Item(124, 14).Value = -119798.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 5251
$ws.Cells.Item(22, 9).Value = 5251
$ws.Cells.Item(22, 11).Value = 5251
$ws.Cells.Item(22, 13).Value = -4956

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 5251
$ws.Cells.Item(27, 9).Value = 5251
$ws.Cells.Item(27, 11).Value = 5251
$ws.Cells.Item(27, 13).Value = -5144

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 4146.2
$ws.Cells.Item(40, 9).Value = 1865.1666
$ws.Cells.Item(40, 10).Value = 5666.8887
$ws.Cells.Item(40, 11).Value = 1865.1666
$ws.Cells.Item(40, 12).Value = 5666.8887
$ws.Cells.Item(40, 13).Value = -1729.1666
$ws.Cells.Item(40, 14).Value = -5938.8887

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(23, 8).Value = 2991.2856
$ws.Cells.Item(23, 9).Value = 688
$ws.Cells.Item(23, 10).Value = 8749.5
$ws.Cells.Item(23, 11).Value = 688
$ws.Cells.Item(23, 12).Value = 8749.5
$ws.Cells.Item(23, 13).Value = -459
$ws.Cells.Item(23, 14).Value = -9207.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 21740498
$ws.Cells.Item(107, 9).Value = 29413320
$ws.Cells.Item(107, 11).Value = 88239960
$ws.Cells.Item(107, 13).Value = -88238040
